# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts.
#
# Updates the DAMSLTag (column I) and DialogAct (column J) values for a
# set of rows on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 13; DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 31; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 32; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 41; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 42; DAMSLTag = "b";  DialogAct = "Acknowledge (Backchannel)" },
    @{ Row = 60; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 74; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 78; DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 79; DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 85; DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 88; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 94; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
)

foreach ($u in $updates) {
    $ws.Range("I$($u.Row)").Value = $u.DAMSLTag
    $ws.Range("J$($u.Row)").Value = $u.DialogAct
}
